# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray empty inline-string cell at C86 ---
$ws.Range("C86").ClearContents()

# --- Ensure the "ano" (D) column cells for the new rows are stored as text ---
$ws.Range("D106:D111").NumberFormat = "@"

# --- Row 106 ---
$ws.Range("A106").Value = 106
$ws.Range("B106").Value = 'CESGRANRIO'
$ws.Range("C106").Value = 'BANCO DO BRASIL'
$ws.Range("D106").Value = '2023'
$ws.Range("E106").Value = '<b>A história do método braile</b>
Ler no escuro. Quem já tentou sabe que é impossível. Mas foi exatamente a isso que um francês chamado Louis Braille dedicou a vida. Nascido em Coupvray, uma pequena aldeia nos arredores de Paris, em 1809, desde cedo ele mostrou muito interesse pelo trabalho do pai. Seus olhos azuis brilhavam da admiração de vê-lo cortar, com extrema perícia, selas e arreios. Pouco depois de completar 3 anos, o menino começou a brincar na selaria do pai, cortando pequenas tiras de couro. Uma tarde, uma sovela, instrumento usado para perfurar o couro, escapou-lhe da mão e atingiu o seu olho esquerdo. O resultado foi uma infecção que, seis meses depois, afetaria também o olho direito. Aos 5 anos, o garoto estava completamente cego.
A tragédia não o impediu, porém, de frequentar a escola por dois anos e de se tornar ainda um aluno brilhante. Por essa razão, ele ganhou uma bolsa de estudos no Instituto Nacional para Jovens Cegos, em Paris, um colégio interno fundado por Valentin Haüy (1745-182. Além do currículo normal, Haüy introduzira um sistema especial de alfabetização, no qual letras de forma impressas em relevo, em papelão, eram reconhecidas pelos contornos. Desde o início do curso, Braille destacou-se como o melhor aluno da turma e logo começou a ajudar os colegas. Em 1821, aos 12 anos, conheceu um método inventado pouco antes por Charles Barbier de La Serre, oficial do Exército francês.
O método Barbier, também chamado escrita noturna, era um código de pontos e traços em relevo impressos também em papelão. Destinava-se a enviar ordens cifradas a sentinelas em postos avançados. Estes decodificariam a mensagem até no escuro. Mas, como a ideia não pegou na tropa, Barbier adaptou o método para a leitura de cegos, com o nome de grafia sonora. O sistema permitia a comunicação entre os cegos, pois com ele era possível escrever, algo que o método de Haüy não possibilitava. O de Barbier era fonético: registrava sons e não letras. Dessa forma, as palavras não podiam ser soletradas. Além disso, o fato de um grande número de sinais ser usado para uma única palavra tornava o sistema muito complicado. Apesar dos inconvenientes, foi adotado como método auxiliar por Haüy.
Pesquisando a fundo a grafia sonora, Braille percebeu suas limitações e pôs-se a aperfeiçoá-la.
Em 1824, seu método estava pronto. Primeiro, eliminou os traços, para evitar erros de leitura: em seguida, criou uma célula de seis pontos, divididos em duas colunas de três pontos cada, que podem ser combinados de 63 maneiras diferentes. A posição dos pontos na célula está ao lado.
Em 1826, aos 17 anos, ainda estudante, Braille começou a dar aulas. Embora seu método fizesse sucesso entre os alunos, não podia ensiná-lo na sala de aula, pois ainda não era reconhecido oficialmente. Por isso, Braille dava aulas do revolucionário sistema escondido no quarto, que logo se transformou numa segunda sala de aula.
O braile é lido passando-se a ponta dos dedos sobre os sinais de relevo. Normalmente se usa a mão direita com um ou mais dedos, conforme a habilidade do leitor, enquanto a mão esquerda procura o início da outra linha. Aplica-se a qualquer língua, sem exceção, e também à estenografia, à música – Braille, por sinal, era ainda exímio pianista – e às notações científicas em geral. A escrita é feita mediante o uso da reglete, também idealizada por Braille: trata-se de uma régua especial, de duas linhas, com uma série de janelas de seis furos cada, correspondentes às células braile.
Louis Braille morreu de tuberculose em 1852, com apenas 43 anos. Temia que seu método desaparecesse com ele, mas, finalmente, em 1854 foi oficializado pelo governo francês. No ano seguinte, foi apresentado ao mundo, na Exposição Internacional de Paris, por ordem do imperador Napoleão III (1808-187, que programou ainda uma série de concertos de piano com ex-alunos de Braille. O sucesso foi imediato, e o sistema se espalhou pelo mundo. Em 1952, o governo francês transferiu os restos mortais de Braille para o Panthéon, em Paris, onde estão sepultados os heróis nacionais.
(ATANES,Silvio.SuperInteressante.Disponívelem:https://super.abril.]com.br/historia/.Acessoem:23out.2022.
Adaptado).
O trecho do parágrafo “Pesquisando a fundo a grafia sonora, Braille percebeu suas limitações e pôs-se a aperfeiçoá- la” pode ser reescrito, sem alterar o sentido que apresenta no texto, como:'
$ws.Range("F106").Value = 'Português'
$ws.Range("G106").Value = 'Emprego de Tempos e Modos'
$ws.Range("H106").Value = 'Médio'
$ws.Range("I106").Value = 'ME'
$ws.Range("J106").Value = 'Para pesquisar a fundo a grafia sonora, Braille percebeu suas limitações e pôs-se a aperfeiçoá-la'
$ws.Range("K106").Value = 'Embora pesquisasse a fundo a grafia sonora, Braille percebeu suas limitações e pôs-se a aperfeiçoá-la'
$ws.Range("L106").Value = 'Quando pesquisava a fundo a grafia sonora, Braille percebeu suas limitações e pôs-se a aperfeiçoá-la'
$ws.Range("M106").Value = 'Apesar de pesquisar a fundo a grafia sonora, Braille percebia suas limitações e punha-se a aperfeiçoá-la'
$ws.Range("N106").Value = 'Se pesquisasse a fundo a grafia sonora, Braille perceberia suas limitações e pôr-se-ia a aperfeiçoá-la'
$ws.Range("O106").Value = 'C'
$ws.Range("P106").Value = 0
$ws.Range("Q106").Value = 0

# --- Row 107 ---
$ws.Range("A107").Value = 107
$ws.Range("B107").Value = 'CESGRANRIO'
$ws.Range("C107").Value = 'ELETRONUCLEAR'
$ws.Range("D107").Value = '2022'
$ws.Range("E107").Value = '<b>Maria José</b>
<i>Paulo Mendes Campos</i>
Faz um ano que Maria José morreu. Era meiga quase sempre, violenta quando necessário. Eu era menino e apanhava de um companheiro maior, quando ela me gritou da sacada se eu não via a pedra que marcava o gol. Dei uma pedrada no outro e acabei com a briga por milagre.
Visitava os miseráveis, internava indigentes enfermos, devotava-se ao alívio de misérias físicas e morais do próximo, estudava o mistério teológico, exigia sempre o mais difícil de si mesma, comungava todos os dias, ingressou na Ordem Terceira de São Francisco. Mas nunca deixou de ter na gaveta o revólver que havia recebido, m enina-e-moça, das mãos do pai, e que empunhou no quintal noturno, perseguindo um ladrão, para espanto de meus cinco anos.
Já perto dos setenta anos, ela explicava para um amigo meu que tinha chegado à humildade da velhice; já não se importava com quem tentasse ofendê-la, mas conservava o revólver para a defesa dos filhos e dos netos.
Tratou-me com a dureza e o carinho que mereciam a rebeldia e o verdor da minha meninice.
Ensinou- me a ler as primeiras sentenças; me falava do Cura d’Ars e nos dois Franciscos, o de Sales e o de Assis; apresentou-me aos contos de Edgar Poe e aos poemas de Baudelaire; dizia-me sorrindo versos de Antônio Nobre que havia decorado quando menina; discutia comigo as ideias finais de Tolstoi; escutava maternalmente meus contos toscos. Quando me desgarrei nos primeiros envolvimentos adolescentes, Maria José, com irônico afeto, me repetia a advertência de Drummond: “Paulo, sossegue, o amor é isso que você está vendo: hoje beija, amanhã não beija, depois de amanhã é domingo e segunda-feira ninguém sabe o que será”.
Logo que me fiz homenzinho, deixou a dureza e se fez minha amiga: nada me perguntava, adivinhava tudo.
Terna e firme, nunca lhe vi a fraqueza da pieguice. Com o gosto espontâneo da qualidade das coisas, renunciou às vaidades mais singelas. Sensível, alegre, aprendeu a encarar o sofrimento de olhos lúcidos. Fiel à disciplina religiosa, compreendia celestialmente as almas que perdiam o rumo. Fé, Esperança e Caridade eram para ela a flecha e o alvo das criaturas.
Tornara-se tão íntima da substância terrestre – a dor – que se fazia difícil para o médico saber o que sentia; acabava dizendo que doía um pouco, por delicadeza.
Capaz de longos jejuns e abstinências, já no final da vida, podia acompanhar um casal amigo a Copacabana, passar do bar da moda ao restaurante diferente, beber dois cafés ou três uísques em santa serenidade e aceitar com alegria o prato exótico.
Gostava das pessoas erradas, consumidas de paixão, admirava São Paulo e Santo Agostinho, acreditava que era preciso se fazer violência para entrar no reino celeste.
Poucas horas antes de morrer, pediu um conhaque e sorriu, destemida e doce, como quem vai partir para o céu. Santificara-se. Deus era o dia e a noite de seu coração, o Pai, a piedade, o fogo do espírito. Perdi quem me amava e perdoava, quem me encomendava à compaixão do Criador e me defendia contra o mundo de revólver na mão.
Disponívelem:https://cronicabrasileira.org.br/cronicas/7173/maria-jose.Acessoem:05fev.2022.
No trecho: “Mas nunca deixou de ter na gaveta o revólver que recebera, menina-e-moça, das mãos do pai, e que empunhou no quintal noturno, perseguindo um ladrão”, (parágrafo 2), a oração destacada pode ser substituída, sem prejuízo de seu significado, por'
$ws.Range("F107").Value = 'Português'
$ws.Range("G107").Value = 'Emprego de Tempos e Modos'
$ws.Range("H107").Value = 'Médio'
$ws.Range("I107").Value = 'ME'
$ws.Range("J107").Value = 'por isso perseguia um ladrão'
$ws.Range("K107").Value = 'enquanto perseguia um ladrão'
$ws.Range("L107").Value = 'embora perseguisse um ladrão'
$ws.Range("M107").Value = 'desde que perseguisse um ladrão'
$ws.Range("N107").Value = 'por mais que perseguisse um ladrão'
$ws.Range("O107").Value = 'B'
$ws.Range("P107").Value = 0
$ws.Range("Q107").Value = 0

# --- Row 108 ---
$ws.Range("A108").Value = 108
$ws.Range("B108").Value = 'CESGRANRIO'
$ws.Range("C108").Value = 'BANCO DO BRASIL'
$ws.Range("D108").Value = '2021'
$ws.Range("E108").Value = 'No trecho “Esse limite <b>poderia </b>ser dado pelo próprio consumidor, se ele assim quiser?” (parágrafo 6), a forma verbal destacada expressa a noção de'
$ws.Range("F108").Value = 'Português'
$ws.Range("G108").Value = 'Emprego de Tempos e Modos'
$ws.Range("H108").Value = 'Médio'
$ws.Range("I108").Value = 'ME'
$ws.Range("J108").Value = 'dever'
$ws.Range("K108").Value = 'certeza'
$ws.Range("L108").Value = 'hipótese'
$ws.Range("M108").Value = 'obrigação'
$ws.Range("N108").Value = 'necessidade'
$ws.Range("O108").Value = 'C'
$ws.Range("P108").Value = 0
$ws.Range("Q108").Value = 0

# --- Row 109 ---
$ws.Range("A109").Value = 109
$ws.Range("B109").Value = 'CESGRANRIO'
$ws.Range("C109").Value = 'BANCO DA AMAZÔNIA'
$ws.Range("D109").Value = '2021'
$ws.Range("E109").Value = 'Em que frase o verbo destacado está flexionado, quanto a número e pessoa, de acordo com a norma-padrão da língua portuguesa?'
$ws.Range("F109").Value = 'Português'
$ws.Range("G109").Value = 'Emprego de Tempos e Modos'
$ws.Range("H109").Value = 'Médio'
$ws.Range("I109").Value = 'ME'
$ws.Range("J109").Value = 'No texto, <u>relacionam</u>-se aos chicles a ideia de eternidade'
$ws.Range("K109").Value = '<u>Referiu</u>-se à eternidade, sem se dar conta, as duas meninas'
$ws.Range("L109").Value = '<u>Enganam</u>-se a respeito da eternidade aqueles que creem nela'
$ws.Range("M109").Value = 'Todos os anos, <u>consome</u>-se muitas balas e chicletes em todo o país'
$ws.Range("N109").Value = 'Em muitas culturas, <u>defendem</u>-se calorosamente a existência da eternidade'
$ws.Range("O109").Value = 'C'
$ws.Range("P109").Value = 0
$ws.Range("Q109").Value = 0

# --- Row 110 ---
$ws.Range("A110").Value = 110
$ws.Range("B110").Value = 'CESGRANRIO'
$ws.Range("C110").Value = 'LIQUIGÁS'
$ws.Range("D110").Value = '2018'
$ws.Range("E110").Value = 'Considere o trecho “Depois vieram as mães e avós doentes.” (l. 8-9).
A frase em que se emprega uma flexão do verbo destacado, de acordo com a norma-padrão da língua portuguesa, é:'
$ws.Range("F110").Value = 'Português'
$ws.Range("G110").Value = 'Emprego de Tempos e Modos'
$ws.Range("H110").Value = 'Médio'
$ws.Range("I110").Value = 'ME'
$ws.Range("J110").Value = 'Não sei o que fazer depois que vinherem as mães e avós doentes'
$ws.Range("K110").Value = 'Depois que as mães e avós doentes virem, faremos alguma coisa'
$ws.Range("L110").Value = 'Depois que eu vim, as mães e avós doentes ficaram curadas'
$ws.Range("M110").Value = 'Depois, as mães e avós doentes tiveram vindo até aqui'
$ws.Range("N110").Value = 'Talvez seja melhor ir depois de vierem as mães e avós doentes'
$ws.Range("O110").Value = 'C'
$ws.Range("P110").Value = 0
$ws.Range("Q110").Value = 0

# --- Row 111 ---
$ws.Range("A111").Value = 111
$ws.Range("B111").Value = 'CESGRANRIO'
$ws.Range("C111").Value = 'TRANSPETRO / SEGUN. OFICIAL'
$ws.Range("D111").Value = '2016'
$ws.Range("E111").Value = 'O distanciamento do autor em relação à história narrada para destacar um ponto de vista seu sobre a temática em foco é marcado pelo uso do verbo <u>ser</u>, no período “<u>É </u>um exercício estranho esse de começar a remoçar um corpo na imaginação, injetar movimento e desejo nos seus músculos, acelerando nele, de novo, a avareza de viver cada instante.” Caso o enunciador queira conferir ao trecho um caráter de possibilidade, a reescritura adequada à norma-padrão e ao contexto empregará o verbo <u>ser </u>da seguinte forma:'
$ws.Range("F111").Value = 'Português'
$ws.Range("G111").Value = 'Emprego de Tempos e Modos'
$ws.Range("H111").Value = 'Médio'
$ws.Range("I111").Value = 'ME'
$ws.Range("J111").Value = 'Fosse'
$ws.Range("K111").Value = 'Seria'
$ws.Range("L111").Value = 'Foi'
$ws.Range("M111").Value = 'Era'
$ws.Range("N111").Value = 'Fora'
$ws.Range("O111").Value = 'B'
$ws.Range("P111").Value = 0
$ws.Range("Q111").Value = 0

